$wb = $excel.ActiveWorkbook

# Updated "想去人数" (F) / "最低票价" (G) figures for matching events on both
# the "展览" sheet and the aggregated "全部类型" sheet.

# Sheet "展览": row -> F (new), G (new). $null means "leave unchanged".
$exhibitionUpdates = @{
    2  = @{ F = $null; G = 65 }
    3  = @{ F = 146;   G = 39 }
    4  = @{ F = 2119;  G = 70 }
    8  = @{ F = 2101;  G = $null }
    10 = @{ F = 10912; G = $null }
    11 = @{ F = 183;   G = $null }
    15 = @{ F = 10709; G = $null }
    16 = @{ F = 430;   G = $null }
    20 = @{ F = 5358;  G = $null }
    22 = @{ F = 3386;  G = $null }
}

# Sheet "全部类型": row -> F (new), G (new).
$allTypesUpdates = @{
    2  = @{ F = $null; G = 65 }
    3  = @{ F = 146;   G = 39 }
    4  = @{ F = 2119;  G = 70 }
    9  = @{ F = 2101;  G = $null }
    13 = @{ F = 10912; G = $null }
    14 = @{ F = 183;   G = $null }
    18 = @{ F = 10709; G = $null }
    19 = @{ F = 430;   G = $null }
    23 = @{ F = 5358;  G = $null }
    25 = @{ F = 3386;  G = $null }
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $vals = $exhibitionUpdates[$row]
    if ($null -ne $vals.F) { $wsExhibition.Range("F$row").Value = $vals.F }
    if ($null -ne $vals.G) { $wsExhibition.Range("G$row").Value = $vals.G }
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $vals = $allTypesUpdates[$row]
    if ($null -ne $vals.F) { $wsAllTypes.Range("F$row").Value = $vals.F }
    if ($null -ne $vals.G) { $wsAllTypes.Range("G$row").Value = $vals.G }
}
